$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 - resale numbers snapshot for 2023-06-03 13:47:07 (Saturday, week 22)
# Columns A-D are textual (date/time/weekday/week-number stored as literal text,
# matching the existing rows' inlineStr convention). Using a leading apostrophe
# via .Formula forces literal-text storage (avoids auto date/number coercion),
# then ClearFormats() strips the quote-prefix style flag so no stray cell style
# is left behind, matching the unstyled look of the other data rows.
$textCells = @{
    "A18" = "2023-06-03"
    "B18" = "13:47:07"
    "C18" = "Saturday"
    "D18" = "22"
}
foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.Formula = "'" + $textCells[$addr]
    $cell.ClearFormats()
}

# Columns E-T are plain numbers.
$ws.Range("E18").Value = 120763
$ws.Range("F18").Value = 133945
$ws.Range("G18").Value = 158821
$ws.Range("H18").Value = 130334
$ws.Range("I18").Value = 174588
$ws.Range("J18").Value = 112471
$ws.Range("K18").Value = 199541
$ws.Range("L18").Value = 218447
$ws.Range("M18").Value = 171690
$ws.Range("N18").Value = 118828
$ws.Range("O18").Value = 38082
$ws.Range("P18").Value = 34825
$ws.Range("Q18").Value = 50146
$ws.Range("R18").Value = -1
$ws.Range("S18").Value = 36712
$ws.Range("T18").Value = -1
